# Correct typo in example data: "Adress Line 1" -> "Address Line 1"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D1")
$cell.Value = "Address Line 1"
$null = $cell.Select()
